$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I7").Value = -0.0165
$ws.Range("J7").Value = -0.0037
$ws.Range("K7").Value = 0.1672
$ws.Range("L7").Value = 0.0786
$ws.Range("M7").Value = 0.0371
$ws.Range("N7").Value = -0.0128
$ws.Range("O7").Value = -0.0628
$ws.Range("P7").Value = 0.0278
$ws.Range("Q7").Value = 0.0256

$ws.Range("I8").Value = -0.4514
$ws.Range("J8").Value = -0.2684
$ws.Range("K8").Value = -0.1513
$ws.Range("L8").Value = -0.0996
$ws.Range("M8").Value = -0.0022
$ws.Range("N8").Value = 0.2975
$ws.Range("O8").Value = 0.3056
$ws.Range("P8").Value = 0.3289
$ws.Range("Q8").Value = 0.3225

$ws.Range("I12").Value = -0.5935
$ws.Range("J12").Value = -0.1412
$ws.Range("K12").Value = -0.015
$ws.Range("L12").Value = 0.0539
$ws.Range("M12").Value = -0.0446
$ws.Range("N12").Value = -0.0601
$ws.Range("O12").Value = -0.0578
$ws.Range("P12").Value = -0.0556
$ws.Range("Q12").Value = -0.0351

$ws.Range("I16").Value = -1.9323
$ws.Range("J16").Value = -1.6849
$ws.Range("K16").Value = -2.403
$ws.Range("L16").Value = -0.8625
$ws.Range("M16").Value = -0.4285
$ws.Range("N16").Value = -0.1597
$ws.Range("O16").Value = -0.3649
$ws.Range("P16").Value = -0.3748
$ws.Range("Q16").Value = -0.0532

$ws.Range("I24").Value = -0.0159
$ws.Range("J24").Value = 0.0745
$ws.Range("K24").Value = -0.0054
$ws.Range("L24").Value = 0.1739
$ws.Range("M24").Value = 0.1421
$ws.Range("N24").Value = 0.1185
$ws.Range("O24").Value = 0.0954
$ws.Range("P24").Value = -0.0667
$ws.Range("Q24").Value = -0.0538

$ws.Range("I35").Value = 0.0071
$ws.Range("J35").Value = 0.0148
$ws.Range("K35").Value = 0.0227
$ws.Range("L35").Value = 0.0295
$ws.Range("M35").Value = 0.0424
$ws.Range("N35").Value = 0.0357
$ws.Range("O35").Value = 0.0288
$ws.Range("P35").Value = 0.0187
$ws.Range("Q35").Value = 0.0048

$ws.Range("I36").Value = 0.236
$ws.Range("J36").Value = 0.1518
$ws.Range("K36").Value = 0.1694
$ws.Range("L36").Value = 0.1514
$ws.Range("M36").Value = 0.1347
$ws.Range("N36").Value = 0.1275
$ws.Range("O36").Value = 0.1296
$ws.Range("P36").Value = 0.0063
$ws.Range("Q36").Value = 0.0168

$ws.Range("I40").Value = 0.1852
$ws.Range("J40").Value = 0.1992
$ws.Range("K40").Value = 0.1946
$ws.Range("L40").Value = 0.0499
$ws.Range("M40").Value = 0.0092
$ws.Range("N40").Value = -0.0089
$ws.Range("O40").Value = -0.0087
$ws.Range("P40").Value = -0.0085
$ws.Range("Q40").Value = -0.0083

$ws.Range("I44").Value = 0.3982
$ws.Range("J44").Value = 0.3253
$ws.Range("K44").Value = 0.3624
$ws.Range("L44").Value = 0.2401
$ws.Range("M44").Value = 0.2028
$ws.Range("N44").Value = 0.1745
$ws.Range("O44").Value = 0.1732
$ws.Range("P44").Value = -0.0016
$ws.Range("Q44").Value = 0.0174

$ws.Range("I52").Value = -0.0391
$ws.Range("J52").Value = -0.0414
$ws.Range("K52").Value = -0.0439
$ws.Range("L52").Value = -0.0147
$ws.Range("M52").Value = -0.011
$ws.Range("N52").Value = -0.0073
$ws.Range("O52").Value = -0.0036
$ws.Range("P52").Value = 0.0032
$ws.Range("Q52").Value = 0.0029
